# Applies the "Stand nach ersten ganzen druchlauf" edit to All_missions.xlsx
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# --- Row 8: category text updates ---
# F8 was "SG/ST" -> now "SG/ST/RL"
$ws.Range("F8").Value = "SG/ST/RL"
# G8 was the old "gesicht muss..." note -> now "Kartendaten"
$ws.Range("G8").Value = "Kartendaten"

# --- Row 12: move the "done" mark from I (Semi) to H (Done) ---
$ws.Range("H12").Value = 1
$ws.Range("I12").ClearContents()

# --- Row 25: same swap as row 12 ---
$ws.Range("H25").Value = 1
$ws.Range("I25").ClearContents()

# --- Row 31: clear the stray "find_face" note in E31, and swap H/I like above ---
$ws.Range("E31").ClearContents()
$ws.Range("H31").Value = 1
$ws.Range("I31").ClearContents()

# --- Row 34: same H/I swap ---
$ws.Range("H34").Value = 1
$ws.Range("I34").ClearContents()

# --- M4 comment box: replace the long CSV-logging text with the shorter
#     "Verwende die vorgegebenen Board-layouts..." note that already exists
#     elsewhere in the workbook ---
$ws.Range("M4").Value = "•`tVerwende die vorgegebenen Board-layouts und QR-Code-PDFs, um deinen Code auszupro-bieren. `n•`tWelche Sensoren von Zumi sind für welche Aufgabe geeignet? `n•`tSpeichern Sie Ihren Code in einem Git-Repository und gewähren Sie den DS-Fachleuten Zu-gang.`n•`tÜbermitteln Sie Ihre Logging-Informationen im CSV-Format an die Dropbox. In diesem Fall ist die Dropbox Ihr GitHub-Repository. Eine Anleitung dazu finden Sie auf Spaces. "

# --- M23 comment box: was empty, now holds the face-detection note ---
$ws.Range("M23").Value = "•`tWenn ein Gesicht erkannt wird, speichern Sie das gescannte Bild als PNG-Datei mit dem Zeit-stempel als Dateinamen"

# --- View state: scroll position & selection, matching the saved workbook state ---
$ws.Activate()
$ws.Range("G17").Select()
$excel.ActiveWindow.ScrollRow = 17
$excel.ActiveWindow.ScrollColumn = 7
$ws.Range("M23:W32").Select()
